# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.827.85'
$ws.Range('E2').Value = '  -4.44%  '

$ws.Range('D3').Value = '2.458.48'
$ws.Range('E3').Value = '  -5.75%  '

$ws.Range('E4').Value = '  +0.20%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '544.08'
$ws.Range('E5').Value = '  -5.30%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.62'
$ws.Range('E6').Value = '  -8.04%  '

$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.611'
$ws.Range('E8').Value = '  -1.78%  '

$ws.Range('D9').Value = '2.453.70'
$ws.Range('E9').Value = '  -5.85%  '

$ws.Range('E10').Value = '  -10.89%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.153'
$ws.Range('E11').Value = '  -1.93%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.34'
$ws.Range('E12').Value = '  -8.66%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.350'
$ws.Range('E13').Value = '  -7.99%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.81'
$ws.Range('E14').Value = '  -8.55%  '

$ws.Range('D15').Value = '2.911.91'
$ws.Range('E15').Value = '  -5.35%  '

$ws.Range('D16').Value = '60.757.61'
$ws.Range('E16').Value = '  -4.38%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000161'
$ws.Range('E17').Value = '  -9.99%  '

$ws.Range('D18').Value = '2.476.00'
$ws.Range('E18').Value = '  -5.78%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.02'
$ws.Range('E19').Value = '  -8.43%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.91'
$ws.Range('E20').Value = '  -9.15%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.14'
$ws.Range('E21').Value = '  -8.78%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '317.74'
$ws.Range('E22').Value = '  -7.48%  '

$ws.Range('E23').Value = '  +0.01%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.15'
$ws.Range('E24').Value = '  -6.52%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.71'
$ws.Range('E25').Value = '  -4.99%  '

$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '2.613.08'
$ws.Range('E26').Value = '  -3.98%  '

$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').Value = '0.0₃0958'
$ws.Range('E27').Value = '  -12.07%  '

$ws.Range('E28').Value = '  +0.31%  '

$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '528.19'
$ws.Range('E29').Value = '  -9.85%  '

$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.23'
$ws.Range('E30').Value = '  -10.27%  '

$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.45'
$ws.Range('E31').Value = '  -8.22%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.57'
$ws.Range('E32').Value = '  -4.26%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.147'
$ws.Range('E33').Value = '  -8.78%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.88'
$ws.Range('E34').Value = '  -8.71%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.54'
$ws.Range('E35').Value = '  -10.86%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.73'
$ws.Range('E36').Value = '  -13.01%  '

$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.26%  '

$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.81'
$ws.Range('E38').Value = '  -10.45%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.372'
$ws.Range('E39').Value = '  -7.86%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.27'
$ws.Range('E40').Value = '  -7.48%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '143.04'
$ws.Range('E41').Value = '  -7.18%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.06%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.67'
$ws.Range('E43').Value = '  -10.70%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.10'
$ws.Range('E44').Value = '  -3.06%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.26'
$ws.Range('E45').Value = '  -9.83%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '146.86'
$ws.Range('E46').Value = '  -6.56%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.54'
$ws.Range('E47').Value = '  -9.34%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.65'
$ws.Range('E48').Value = '  -13.32%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0527'
$ws.Range('E49').Value = '  -10.49%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.582'
$ws.Range('E50').Value = '  -7.60%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0935'
$ws.Range('E51').Value = '  -6.68%  '
